$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update model names in column A (rows 2-26) per the permutation
$ws.Range("A2").Value = "model_4_6_0"
$ws.Range("A3").Value = "model_4_6_22"
$ws.Range("A4").Value = "model_4_6_21"
$ws.Range("A5").Value = "model_4_6_20"
$ws.Range("A6").Value = "model_4_6_19"
$ws.Range("A7").Value = "model_4_6_18"
$ws.Range("A8").Value = "model_4_6_17"
$ws.Range("A9").Value = "model_4_6_16"
$ws.Range("A10").Value = "model_4_6_15"
$ws.Range("A11").Value = "model_4_6_14"
$ws.Range("A12").Value = "model_4_6_13"
$ws.Range("A13").Value = "model_4_6_23"
$ws.Range("A14").Value = "model_4_6_12"
$ws.Range("A15").Value = "model_4_6_10"
$ws.Range("A16").Value = "model_4_6_9"
$ws.Range("A17").Value = "model_4_6_8"
$ws.Range("A18").Value = "model_4_6_7"
$ws.Range("A19").Value = "model_4_6_6"
$ws.Range("A20").Value = "model_4_6_5"
$ws.Range("A21").Value = "model_4_6_4"
$ws.Range("A22").Value = "model_4_6_3"
$ws.Range("A23").Value = "model_4_6_2"
$ws.Range("A24").Value = "model_4_6_1"
$ws.Range("A25").Value = "model_4_6_11"
$ws.Range("A26").Value = "model_4_6_24"

# Update metric columns B:Q for all rows 2-26 to the same new constant values
$rowValues = @{
    "B" = 0.445980542904741
    "C" = -0.05067851974480586
    "D" = 0.7720579009814785
    "E" = 0.555840492537803
    "F" = 0.7483261297087673
    "G" = 0.3288898821261522
    "H" = 0.6237281562693088
    "I" = 0.2715281851569437
    "J" = 0.07442302943515401
    "K" = 0.1729756072960488
    "L" = 0.279043433096329
    "M" = 0.5734892170966707
    "N" = 0.05025235926527027
    "O" = 0.5979038403548141
    "P" = 22.22406457805329
    "Q" = 34.41282282673529
}

for ($row = 2; $row -le 26; $row++) {
    $ws.Range("B$row").Value = $rowValues["B"]
    $ws.Range("C$row").Value = $rowValues["C"]
    $ws.Range("D$row").Value = $rowValues["D"]
    $ws.Range("E$row").Value = $rowValues["E"]
    $ws.Range("F$row").Value = $rowValues["F"]
    $ws.Range("G$row").Value = $rowValues["G"]
    $ws.Range("H$row").Value = $rowValues["H"]
    $ws.Range("I$row").Value = $rowValues["I"]
    $ws.Range("J$row").Value = $rowValues["J"]
    $ws.Range("K$row").Value = $rowValues["K"]
    $ws.Range("L$row").Value = $rowValues["L"]
    $ws.Range("M$row").Value = $rowValues["M"]
    $ws.Range("N$row").Value = $rowValues["N"]
    $ws.Range("O$row").Value = $rowValues["O"]
    $ws.Range("P$row").Value = $rowValues["P"]
    $ws.Range("Q$row").Value = $rowValues["Q"]
}
